$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Create row 12 (new row at the bottom), copying formats from row 11 ---
$ws.Range("A11:G11").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Give the new dates on row 8 (Config CI/CD) a date number format like the other date cells ---
$ws.Range("C8:D8").NumberFormat = "d-mmm-yy"

# --- Row 6: new "Unit test" task (keeps the existing dates 9-Sep-19 / 13-Sep-19) ---
$ws.Range("B6").Value = "Unit test"
$ws.Range("G6").ClearContents()

# --- Row 7: "Deploy docker local" moves here, with the dates formerly on Config CI/CD ---
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Deploy docker local"
$ws.Range("C7").Value = 43724
$ws.Range("D7").Value = 43728

# --- Row 8: "Config CI/CD" moves here with brand-new dates ---
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = " Config CI/CD"
$ws.Range("C8").Value = 43731
$ws.Range("D8").Value = 43735

# --- Row 9: "Register Azure Account" ---
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "Register Azure Account"

# --- Row 10: "Add Docker Registry in Azure and config Auto CI/CD" ---
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Add Docker Registry in Azure and config Auto CI/CD"

# --- Row 11: "Add Wiki page" ---
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Add Wiki page"
$ws.Range("G11").Value = "Not started"

# --- Row 12 (new): "Send to line manager" ---
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "Send to line manager"
$ws.Range("G12").Value = "Not started"

# --- Selection / active cell like the authored workbook ---
$ws.Range("D9").Select() | Out-Null
